# Fruta / hortaliza, semanal
# Insert two new weekly rows (Espárragos, Mapocho Venta Directa de Santiago,
# "Sin especificar" variety, Primera/Segunda quality, 2021-11-23) above the
# existing row 54, shifting all subsequent data rows (old 54-71) down to
# 56-73, and populate the two new rows with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 54 - shifts existing rows 54:71 down to 56:73
# and copies formatting (incl. the date-style on column D) from the row above.
$ws.Range("54:55").Insert()

# New row 54: Espárragos, Sin especificar, Primera
$ws.Range("A54").Value2 = 12
$ws.Range("B54").Value2 = "Mapocho Venta Directa de Santiago"
$ws.Range("C54").Value2 = "Metropolitana"
$ws.Range("D54").Value2 = 44523
$ws.Range("E54").Value2 = 13
$ws.Range("F54").Value2 = 300000000
$ws.Range("G54").Value2 = "Espárragos"
$ws.Range("H54").Value2 = "Sin especificar"
$ws.Range("I54").Value2 = "Primera"
$ws.Range("J54").Value2 = 480
$ws.Range("K54").Value2 = 1000
$ws.Range("L54").Value2 = 1000
$ws.Range("M54").Value2 = 1000
$ws.Range("N54").Value2 = '$/kilo'
$ws.Range("O54").Value2 = "Provincia de Linares"
$ws.Range("P54").Value2 = 1000
$ws.Range("Q54").Value2 = 1
$ws.Range("R54").Value2 = "Hortaliza"

# New row 55: Espárragos, Sin especificar, Segunda
$ws.Range("A55").Value2 = 12
$ws.Range("B55").Value2 = "Mapocho Venta Directa de Santiago"
$ws.Range("C55").Value2 = "Metropolitana"
$ws.Range("D55").Value2 = 44523
$ws.Range("E55").Value2 = 13
$ws.Range("F55").Value2 = 300000000
$ws.Range("G55").Value2 = "Espárragos"
$ws.Range("H55").Value2 = "Sin especificar"
$ws.Range("I55").Value2 = "Segunda"
$ws.Range("J55").Value2 = 450
$ws.Range("K55").Value2 = 800
$ws.Range("L55").Value2 = 800
$ws.Range("M55").Value2 = 800
$ws.Range("N55").Value2 = '$/kilo'
$ws.Range("O55").Value2 = "Provincia de Linares"
$ws.Range("P55").Value2 = 800
$ws.Range("Q55").Value2 = 1
$ws.Range("R55").Value2 = "Hortaliza"
